$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 106, pushing the old row 106 (blank separator)
# and the summary rows (107-109) down by one. Excel copies formatting from the
# row above (105) onto the newly inserted row, which matches the target
# styling for columns D/E/F/G.
$ws.Range("A106").EntireRow.Insert()

# Fill in the new data row 106 (a work interval on 2014-03-31).
$ws.Range("A106").Value = 2014
$ws.Range("B106").Value = 3
$ws.Range("C106").Value = 31
$ws.Range("D106").Value = 0.73958333333333337
$ws.Range("E106").Value = 0.79166666666666663

# Time-spent / hours-spent formulas for the new row, matching the pattern
# used throughout column F/G.
$ws.Range("F106").Formula = "=(E106-D106)*24*60"
$ws.Range("G106").Formula = "=F106/60"

# Shrink the end time of the previous entry (row 105); dependent formulas
# (F105, G105) and the downstream sums recalculate automatically.
$ws.Range("E105").Value = 0.64930555555555558

# Restore the originally selected cell on the sheet.
$ws.Range("E107").Select()
